$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $text) {
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $text
}

Set-CellText 1 1 "11÷4=2, 3"
Set-CellText 1 2 "34÷4=8, 2"
Set-CellText 1 3 "76÷2=38, 0"
Set-CellText 1 4 "64÷7=9, 1"
Set-CellText 1 5 "10÷5=2, 0"

Set-CellText 5 1 "76÷6=12, 4"
Set-CellText 5 2 "48÷8=6, 0"
Set-CellText 5 3 "98÷4=24, 2"
Set-CellText 5 4 "89÷4=22, 1"
Set-CellText 5 5 "83÷4=20, 3"

Set-CellText 9 1 "23÷9=2, 5"
Set-CellText 9 2 "78÷2=39, 0"
Set-CellText 9 3 "83÷8=10, 3"
Set-CellText 9 4 "52÷6=8, 4"
Set-CellText 9 5 "83÷8=10, 3"

Set-CellText 13 1 "81÷6=13, 3"
Set-CellText 13 2 "91÷6=15, 1"
Set-CellText 13 3 "43÷7=6, 1"
Set-CellText 13 4 "45÷3=15, 0"
Set-CellText 13 5 "19÷8=2, 3"

Set-CellText 17 1 "50÷2=25, 0"
Set-CellText 17 2 "25÷5=5, 0"
Set-CellText 17 3 "44÷9=4, 8"
Set-CellText 17 4 "29÷9=3, 2"
Set-CellText 17 5 "70÷7=10, 0"
